$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All target cells are plain text
# (inline strings) in the source workbook, e.g. "29.391.81" or "  -0.16%  ",
# so we force a text number format before assigning the value and then clear
# the style back to Normal so no stray formatting/style id is left behind.
$updates = [ordered]@{
    'D2' = '29.393.04'
    'E2' = '  -0.12%  '
    'D3' = '1.844.13'
    'E3' = '  -0.25%  '
    'D4' = '0.9992'
    'E4' = '  +0.02%  '
    'D5' = '238.95'
    'E5' = '  -0.80%  '
    'D6' = '0.6313'
    'E6' = '  -0.45%  '
    'E7' = '  +0.02%  '
    'D8' = '0.07526'
    'E8' = '  -0.47%  '
    'D9' = '0.2927'
    'E9' = '  -1.48%  '
    'D10' = '24.54'
    'E10' = '  -0.46%  '
    'D11' = '0.07714'
    'E11' = '  -0.01%  '
    'D12' = '1.843.25'
    'E12' = '  -7.14%  '
    'E13' = '  +0.17%  '
    'D14' = '0.6794'
    'E14' = '  -1.01%  '
    'D15' = '0.00001043'
    'E15' = '  +5.06%  '
    'D16' = '83.31'
    'E16' = '  +0.37%  '
    'D17' = '2.089.73'
    'E17' = '  -7.71%  '
    'D18' = '6.174'
    'E18' = '  -0.41%  '
    'D19' = '29.427.34'
    'E19' = '  -0.15%  '
    'D20' = '228.49'
    'E20' = '  -1.48%  '
    'E21' = '  -0.73%  '
    'E22' = '  +0.04%  '
    'D23' = '7.459'
    'E23' = '  -1.85%  '
    'E24' = '  +0.05%  '
    'D25' = '157.01'
    'E25' = '  +0.91%  '
    'D26' = '0.1394'
    'E26' = '  +0.52%  '
    'D27' = '8.353'
    'E27' = '  -0.88%  '
    'D28' = '17.59'
    'E28' = '  -0.56%  '
    'D29' = '1.456'
    'E29' = '  -0.83%  '
    'D30' = '1.281'
    'E30' = '  +1.81%  '
    'E31' = '  -3.06%  '
    'D32' = '4.103'
    'E32' = '  -0.90%  '
    'D33' = '4.020'
    'E33' = '  -0.01%  '
    'D34' = '1.845'
    'E34' = '  -0.76%  '
    'D35' = '1.157'
    'E35' = '  -0.03%  '
    'D36' = '0.7120'
    'E36' = '  -0.67%  '
    'D37' = '2.591'
    'E37' = '  +0.00%  '
    'D38' = '1.246.00'
    'E38' = '  -0.48%  '
    'D39' = '0.01811'
    'E39' = '  +0.32%  '
    'E40' = '  -1.01%  '
    'D41' = '6.357'
    'E41' = '  +4.34%  '
    'D42' = '0.9026'
    'E42' = '  +0.20%  '
    'E43' = '  +0.02%  '
    'D44' = '101.58'
    'E44' = '  -0.09%  '
    'D45' = '65.73'
    'E45' = '  -1.99%  '
    'B46' = 'BabyDogeCoin'
    'C46' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D46' = '0.00000000118'
    'E46' = '  +0.33%  '
    'B47' = 'Aptos'
    'C47' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D47' = '7.102'
    'E47' = '  -1.44%  '
    'E48' = '  -0.70%  '
    'D49' = '8.963'
    'E49' = '  -1.90%  '
    'D50' = '1.671'
    'E50' = '  -0.98%  '
    'D51' = '0.1122'
    'E51' = '  -0.37%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
